$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.155.36"
$ws.Range("E2").Value = "  +1.46%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.905.85"
$ws.Range("E3").Value = "  +1.84%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.003"
$ws.Range("E4").Value = "  +0.03%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "327.10"
$ws.Range("E5").Value = "  +0.73%  "

$ws.Range("E6").Value = "  +0.03%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4637"
$ws.Range("E7").Value = "  +0.72%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3942"
$ws.Range("E8").Value = "  +2.06%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "46.62"
$ws.Range("E9").Value = "  +0.86%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07969"
$ws.Range("E10").Value = "  +1.38%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.9992"
$ws.Range("E11").Value = "  +1.85%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "22.22"
$ws.Range("E12").Value = "  +1.83%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.947.94"
$ws.Range("E13").Value = "  +3.41%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.106"
$ws.Range("E14").Value = "  +1.51%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.765"
$ws.Range("E15").Value = "  +1.04%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.06960"
$ws.Range("E16").Value = "  -0.05%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "88.50"
$ws.Range("E17").Value = "  +0.07%  "

$ws.Range("E18").Value = "  +0.02%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.00001006"
$ws.Range("E19").Value = "  +0.38%  "

$ws.Range("E20").Value = "  +2.24%  "

$ws.Range("E21").Value = "  +0.06%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "29.191.07"
$ws.Range("E22").Value = "  +1.57%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.366"
$ws.Range("E23").Value = "  +1.76%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.09"
$ws.Range("E24").Value = "  +0.13%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.159.50"
$ws.Range("E25").Value = "  +2.62%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.051"
$ws.Range("E26").Value = "  -2.31%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "156.72"
$ws.Range("E27").Value = "  +2.79%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "19.48"
$ws.Range("E28").Value = "  +0.97%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.987"
$ws.Range("E29").Value = "  +1.63%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.991"
$ws.Range("E30").Value = "  +0.13%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "118.93"
$ws.Range("E31").Value = "  -0.18%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09394"
$ws.Range("E32").Value = "  +0.73%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.9199"
$ws.Range("E33").Value = "  +0.36%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.343"
$ws.Range("E34").Value = "  +0.85%  "

$ws.Range("E35").Value = "  +1.01%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.264"
$ws.Range("E36").Value = "  -1.82%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.05813"
$ws.Range("E37").Value = "  +0.52%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.183"
$ws.Range("E38").Value = "  +3.16%  "

$ws.Range("E39").Value = "  +1.24%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "7.990"
$ws.Range("E40").Value = "  +4.41%  "

$ws.Range("E41").Value = "  +2.11%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1802"
$ws.Range("E42").Value = "  +1.27%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "9.965"
$ws.Range("E43").Value = "  +1.95%  "

$ws.Range("E44").Value = "  +2.69%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.5417"
$ws.Range("E45").Value = "  +2.37%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.198"
$ws.Range("E46").Value = "  +3.32%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.07087"
$ws.Range("E47").Value = "  -1.90%  "

$ws.Range("E48").Value = "  +2.14%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.566"
$ws.Range("E49").Value = "  +6.52%  "

$ws.Range("E50").Value = "  -0.72%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.046"
$ws.Range("E51").Value = "  -6.80%  "
